$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add precio_num / fecha_dia columns ---
$ws.Range("E1").Value = "precio_num"
$ws.Range("E1").Style = $ws.Range("A1").Style
$ws.Range("F1").Value = "fecha_dia"
$ws.Range("F1").Style = $ws.Range("A1").Style

# --- Set the date format used for column F (creates numFmt 166 lowercase,
#     then 167 uppercase -- the uppercase one is what ends up applied) ---
$ws.Range("F2").NumberFormat = "yyyy-mm-dd"
$ws.Range("F2").NumberFormat = "YYYY-MM-DD"

# --- Row 2..107: backfill precio_num (E) and fecha_dia (F) for existing rows ---
For ($r = 2; $r -le 107; $r++) {
    $priceText = $ws.Cells.Item($r, 4).Value2
    $priceNum = [double]($priceText -replace "€","" -replace ",",".")
    $ws.Cells.Item($r, 5).Value = $priceNum
    $dateSerial = [math]::Floor($ws.Cells.Item($r, 1).Value2)
    $ws.Cells.Item($r, 6).Value = $dateSerial
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD"
}

# --- Row 107: the timestamp in A107 was re-entered with a slightly different
#     fractional value ---
$ws.Range("A107").Value = 45964.36351739583

# --- New rows 108..134: freshly captured price entries ---

$ws.Cells.Item(108, 1).Value = 45966.36925192129
$ws.Cells.Item(108, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(108, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(108, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(108, 4).Value = "53,10€"
$ws.Cells.Item(108, 5).Value = 53.1
$ws.Cells.Item(108, 6).Value = 45966
$ws.Cells.Item(108, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(109, 1).Value = 45966.3976755324
$ws.Cells.Item(109, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(109, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(109, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(109, 4).Value = "53,10€"
$ws.Cells.Item(109, 5).Value = 53.1
$ws.Cells.Item(109, 6).Value = 45966
$ws.Cells.Item(109, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(110, 1).Value = 45966.45849217592
$ws.Cells.Item(110, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(110, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(110, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(110, 4).Value = "53,10€"
$ws.Cells.Item(110, 5).Value = 53.1
$ws.Cells.Item(110, 6).Value = 45966
$ws.Cells.Item(110, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(111, 1).Value = 45967.45869030093
$ws.Cells.Item(111, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(111, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(111, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(111, 4).Value = "53,10€"
$ws.Cells.Item(111, 5).Value = 53.1
$ws.Cells.Item(111, 6).Value = 45967
$ws.Cells.Item(111, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(112, 1).Value = 45968.45859172453
$ws.Cells.Item(112, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(112, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(112, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(112, 4).Value = "53,10€"
$ws.Cells.Item(112, 5).Value = 53.1
$ws.Cells.Item(112, 6).Value = 45968
$ws.Cells.Item(112, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(113, 1).Value = 45969.45855380787
$ws.Cells.Item(113, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(113, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(113, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(113, 4).Value = "53,10€"
$ws.Cells.Item(113, 5).Value = 53.1
$ws.Cells.Item(113, 6).Value = 45969
$ws.Cells.Item(113, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(114, 1).Value = 45970.45855238426
$ws.Cells.Item(114, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(114, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(114, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(114, 4).Value = "53,10€"
$ws.Cells.Item(114, 5).Value = 53.1
$ws.Cells.Item(114, 6).Value = 45970
$ws.Cells.Item(114, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(115, 1).Value = 45971.45856265046
$ws.Cells.Item(115, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(115, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(115, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(115, 4).Value = "53,10€"
$ws.Cells.Item(115, 5).Value = 53.1
$ws.Cells.Item(115, 6).Value = 45971
$ws.Cells.Item(115, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(116, 1).Value = 45972.45856766203
$ws.Cells.Item(116, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(116, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(116, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(116, 4).Value = "44,75€"
$ws.Cells.Item(116, 5).Value = 44.75
$ws.Cells.Item(116, 6).Value = 45972
$ws.Cells.Item(116, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(117, 1).Value = 45973.4585937037
$ws.Cells.Item(117, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(117, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(117, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(117, 4).Value = "53,10€"
$ws.Cells.Item(117, 5).Value = 53.1
$ws.Cells.Item(117, 6).Value = 45973
$ws.Cells.Item(117, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(118, 1).Value = 45974.45856204861
$ws.Cells.Item(118, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(118, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(118, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(118, 4).Value = "53,10€"
$ws.Cells.Item(118, 5).Value = 53.1
$ws.Cells.Item(118, 6).Value = 45974
$ws.Cells.Item(118, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(119, 1).Value = 45978.4585547801
$ws.Cells.Item(119, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(119, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(119, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(119, 4).Value = "53,10€"
$ws.Cells.Item(119, 5).Value = 53.1
$ws.Cells.Item(119, 6).Value = 45978
$ws.Cells.Item(119, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(120, 1).Value = 45979.45854807871
$ws.Cells.Item(120, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(120, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(120, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(120, 4).Value = "53,10€"
$ws.Cells.Item(120, 5).Value = 53.1
$ws.Cells.Item(120, 6).Value = 45979
$ws.Cells.Item(120, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(121, 1).Value = 45980.45854332176
$ws.Cells.Item(121, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(121, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(121, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(121, 4).Value = "53,10€"
$ws.Cells.Item(121, 5).Value = 53.1
$ws.Cells.Item(121, 6).Value = 45980
$ws.Cells.Item(121, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(122, 1).Value = 45981.45853957176
$ws.Cells.Item(122, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(122, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(122, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(122, 4).Value = "53,10€"
$ws.Cells.Item(122, 5).Value = 53.1
$ws.Cells.Item(122, 6).Value = 45981
$ws.Cells.Item(122, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(123, 1).Value = 45985.45857799768
$ws.Cells.Item(123, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(123, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(123, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(123, 4).Value = "46,10€"
$ws.Cells.Item(123, 5).Value = 46.1
$ws.Cells.Item(123, 6).Value = 45985
$ws.Cells.Item(123, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(124, 1).Value = 45986.39181813657
$ws.Cells.Item(124, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(124, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(124, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(124, 4).Value = "46,10€"
$ws.Cells.Item(124, 5).Value = 46.1
$ws.Cells.Item(124, 6).Value = 45986
$ws.Cells.Item(124, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(125, 1).Value = 45986.40511755787
$ws.Cells.Item(125, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(125, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(125, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(125, 4).Value = "46,10€"
$ws.Cells.Item(125, 5).Value = 46.1
$ws.Cells.Item(125, 6).Value = 45986
$ws.Cells.Item(125, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(126, 1).Value = 45986.40574495371
$ws.Cells.Item(126, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(126, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(126, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(126, 4).Value = "46,10€"
$ws.Cells.Item(126, 5).Value = 46.1
$ws.Cells.Item(126, 6).Value = 45986
$ws.Cells.Item(126, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(127, 1).Value = 45986.40638818287
$ws.Cells.Item(127, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(127, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(127, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(127, 4).Value = "46,10€"
$ws.Cells.Item(127, 5).Value = 46.1
$ws.Cells.Item(127, 6).Value = 45986
$ws.Cells.Item(127, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(128, 1).Value = 45986.40847553241
$ws.Cells.Item(128, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(128, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(128, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(128, 4).Value = "46,10€"
$ws.Cells.Item(128, 5).Value = 46.1
$ws.Cells.Item(128, 6).Value = 45986
$ws.Cells.Item(128, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(129, 1).Value = 45986.41846228009
$ws.Cells.Item(129, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(129, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(129, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(129, 4).Value = "46,10€"
$ws.Cells.Item(129, 5).Value = 46.1
$ws.Cells.Item(129, 6).Value = 45986
$ws.Cells.Item(129, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(130, 1).Value = 45986.41970927083
$ws.Cells.Item(130, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(130, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(130, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(130, 4).Value = "46,10€"
$ws.Cells.Item(130, 5).Value = 46.1
$ws.Cells.Item(130, 6).Value = 45986
$ws.Cells.Item(130, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(131, 1).Value = 45986.42277232639
$ws.Cells.Item(131, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(131, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(131, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(131, 4).Value = "46,10€"
$ws.Cells.Item(131, 5).Value = 46.1
$ws.Cells.Item(131, 6).Value = 45986
$ws.Cells.Item(131, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(132, 1).Value = 45986.42968309027
$ws.Cells.Item(132, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(132, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(132, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(132, 4).Value = "46,10€"
$ws.Cells.Item(132, 5).Value = 46.1
$ws.Cells.Item(132, 6).Value = 45986
$ws.Cells.Item(132, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(133, 1).Value = 45986.43294613426
$ws.Cells.Item(133, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(133, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(133, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(133, 4).Value = "46,10€"
$ws.Cells.Item(133, 5).Value = 46.1
$ws.Cells.Item(133, 6).Value = 45986
$ws.Cells.Item(133, 6).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(134, 1).Value = 45986.43557956163
$ws.Cells.Item(134, 1).Style = $ws.Range("A2").Style
$ws.Cells.Item(134, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(134, 3).Value = "Pack (5x500g)"
$ws.Cells.Item(134, 4).Value = "46,10€"
# Row 134 is the newest entry -- precio_num / fecha_dia not back-filled yet
$ws.Cells.Item(134, 5).Value = ""
$ws.Cells.Item(134, 6).Value = ""
